$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'76.486.09"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.75%  "

$ws.Range("D3").Value = "'2.959.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.71%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("E5").Value = "  +1.82%  "

$ws.Range("D6").Value = "'597.06"
$ws.Range("D6").Style = "Normal"

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("E9").Value = "  +5.78%  "

$ws.Range("D10").Value = "'2.959.12"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.66%  "

$ws.Range("E11").Value = "  +12.05%  "

$ws.Range("E12").Value = "  +0.58%  "

$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "'3.501.97"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.85%  "

$ws.Range("B14").Value = "Toncoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D14").Value = "'4.92"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.72%  "

$ws.Range("D15").Value = "'28.56"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.90%  "

$ws.Range("D16").Value = "'76.380.36"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.81%  "

$ws.Range("D17").Value = "'0.0000190"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.65%  "

$ws.Range("D18").Value = "'2.959.71"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.74%  "

$ws.Range("D19").Value = "'13.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +9.19%  "

$ws.Range("D20").Value = "'8.76"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.88%  "

$ws.Range("D21").Value = "'378.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.32%  "

$ws.Range("E22").Value = "  +0.04%  "

$ws.Range("E23").Value = "  +4.73%  "

$ws.Range("D24").Value = "'72.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.51%  "

$ws.Range("E25").Value = "  -0.16%  "

$ws.Range("E26").Value = "  +2.22%  "

$ws.Range("D27").Value = "'4.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.26%  "

$ws.Range("D28").Value = "'9.74"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.22%  "

$ws.Range("E29").Value = "  +0.66%  "

$ws.Range("E30").Value = "  +0.15%  "

$ws.Range("D31").Value = "'8.59"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +11.03%  "

$ws.Range("D32").Value = "'1.39"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.67%  "

$ws.Range("D33").Value = "'496.77"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.75%  "

$ws.Range("E34").Value = "  +1.38%  "

$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.12%  "

$ws.Range("D36").Value = "'165.36"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.32%  "

$ws.Range("D37").Value = "'20.38"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.16%  "

$ws.Range("E38").Value = "  +14.47%  "

$ws.Range("D39").Value = "'0.109"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +19.48%  "

$ws.Range("E40").Value = "  +1.35%  "

$ws.Range("E41").Value = "  -1.64%  "

$ws.Range("E42").Value = "  +0.00%  "

$ws.Range("D43").Value = "'180.69"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.51%  "

$ws.Range("D44").Value = "'4.95"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.58%  "

$ws.Range("E45").Value = "  -1.12%  "

$ws.Range("D46").Value = "'39.91"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.24%  "

$ws.Range("E47").Value = "  -1.93%  "

$ws.Range("D48").Value = "'0.591"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.43%  "

$ws.Range("D49").Value = "'3.92"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.70%  "

$ws.Range("D50").Value = "'2.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.80%  "

$ws.Range("D51").Value = "'0.669"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.44%  "
